# SCD0175_Penambahan Leads dari Store ke Cart.xlsx
# Update the "nama" leads sample values in column L (rows 2-5) and move the
# active selection to L6, as per commit:
# "Scripting SCD0179 - Validasi Field Report PHR Pada Searching
# Portal/Action0/ObjectRepository and SCD0180 - Sales Mengakses Menu
# Report - Menu Product Holding Ratio - Report"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace sample "nama" values used by the leads-to-cart test rows.
$ws.Range("L2").Value = "nadia 2"
$ws.Range("L3").Value = "dewi 8"
$ws.Range("L4").Value = "tyas"
$ws.Range("L5").Value = "bnimf"

# Move/save the active cell selection to L6 (matches saved sheet view state).
$ws.Range("L6").Select()
